# Add a new "Save" column (H) to the s_vals sheet, matching the existing
# header formatting used by the other header cells (e.g. G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell text
$ws.Range("H1").Value = "Save"

# Copy the formatting (style) from the neighboring header cell G1 onto H1
# so it reuses the same cell style rather than creating a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data cell for row 2
$ws.Range("H2").Value = 1
